# Small change to make a test:
# Append a new "major" block (Other major) to the "Esempi" sheet, mirroring
# the existing "Business Administration" block that precedes it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Esempi")

# New rows to append after the existing data (rows 1-26).
$newRows = @(
    @("Other major", "MA101", "Something else", 2, $null),
    @("ADD 1", 1, "CS", 110, 110),
    @("ADD 2", 1, "MA", 197, 197),
    @("COR 1", 1, "MA", 209, 209),
    @("COR 2", 1, "BUS ", 301, 301),
    @("ELC 1", 3, "FIN", 300, 1000)
)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $val = $values[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

# Match the saved selection state from the edit (scroll position itself
# isn't persisted through this runtime's writer, but the active selection is).
$ws.Activate()
$ws.Range("D34").Select()
